$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column C of outcome measurements (Pre Experimental Phase values)
$ws.Range("C2").Value = "Not stressful"
$ws.Range("C3").Value = "Moderately stressful"
$ws.Range("C4").Value = "A little stressful"
$ws.Range("C5").Value = "Not stressful"
$ws.Range("C6").Value = "A little stressful"
$ws.Range("C7").Value = "A little stressful"

# Update selection to reflect the active cell after data entry
$ws.Range("C8").Select()
